{"js": "// Gr\u00e1fica 2 fix: justify the Normal style and bold the \"Heading 2\" /\n// \"Header 2\" title styles (paragraph style + their linked character\n// styles), matching the template restyle described in the commit.\n\nconst styles = context.document.getStyles();\n\n// 1) Normal: paragraphs based on Normal should be fully justified.\nconst normal = styles.getByNameOrNullObject(\"Normal\");\n\n// 2) \"heading 2\" (styleId Ttulo2) + its linked character style\n//    \"T\u00edtulo 2 Car\" (styleId Ttulo2Car): make the heading text bold.\nconst heading2 = styles.getByNameOrNullObject(\"Heading 2\");\nconst heading2Char = styles.getByNameOrNullObject(\"T\u00edtulo 2 Car\");\n\n// 3) Custom \"Header 2\" style's linked character style \"Header 2 Car\"\n//    (styleId Header2Car): make it bold too.\nconst header2Char = styles.getByNameOrNullObject(\"Header 2 Car\");\n\nawait context.sync();\n\nif (!normal.isNullObject) {\n  normal.paragraphFormat.alignment = Word.Alignment.justified;\n}\n\nif (!heading2.isNullObject) {\n  heading2.font.bold = true;\n}\n\nif (!heading2Char.isNullObject) {\n  heading2Char.font.bold = true;\n}\n\nif (!header2Char.isNullObject) {\n  header2Char.font.bold = true;\n}\n\nawait context.sync();\n", "ps1": "# Gr\u00e1fica 2 fix: justify the Normal style and bold the \"Heading 2\" /\n# \"Header 2\" title styles (paragraph style + their linked character\n# styles), matching the template restyle described in the commit.\n\n$d = $word.ActiveDocument\n\n# 1) Normal: paragraphs based on Normal should be fully justified.\n$normal = $d.Styles(\"Normal\")\n$normal.ParagraphFormat.Alignment = 3  # wdAlignParagraphJustify\n\n# 2) \"heading 2\" (styleId Ttulo2) + its linked character style\n#    \"T\u00edtulo 2 Car\" (styleId Ttulo2Car): make the heading text bold.\n$heading2 = $d.Styles(\"Heading 2\")\n$heading2.Font.Bold = $true\n\n$heading2Char = $d.Styles(\"T\u00edtulo 2 Car\")\n$heading2Char.Font.Bold = $true\n\n# 3) Custom \"Header 2\" style's linked character style \"Header 2 Car\"\n#    (styleId Header2Car): make it bold too.\n$header2Char = $d.Styles(\"Header 2 Car\")\n$header2Char.Font.Bold = $true\n"}
